$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '87.975.04'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.78%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.103.32'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.62%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.08'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '634.38'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.75%  '
$ws.Range("E7").Value = '  -0.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.782'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +13.37%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.099.30'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.79%  '
$ws.Range("E11").Value = '  -3.20%  '
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000247'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.35'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.928.55'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.676.71'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.57%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.97'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.117.71'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.36%  '
$ws.Range("E19").Value = '  +2.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000217'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +16.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.15'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '419.82'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.36'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.87'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.41'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.80'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +9.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.43'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.31%  '
$ws.Range("E31").Value = '  -8.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.00'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.12'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '499.13'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.56%  '
$ws.Range("E35").Value = '  +13.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.84'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.26'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.09'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.18'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.362'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.83'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.31%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.133'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +7.58%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '145.76'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.56'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0652'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +11.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '161.27'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -6.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.711'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("E51").Value = '  -5.32%  '
